$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2000-09")
$ws.Activate()

# Insert a new row at row 7 (pushes the old row 7 blank-formatting row and
# row 8 header row down to 8 and 9), inheriting formats from the row above.
$ws.Rows("7:7").Insert()

# --- Fill in the new data row (row 7) ---
$ws.Range("A7").Value = "CW3M"
$ws.Range("B7").Value = "Baseline_2000-09_C366"
$ws.Range("C7").Value = "2000-09"

$ws.Range("D7").Value = 898.93680410000002
$ws.Range("E7").Value = 1790.8486085
$ws.Range("F7").Value = 1.0508375999999999
$ws.Range("G7").Value = 270.46964409999998
$ws.Range("H7").Value = 9.3795183999999985
$ws.Range("I7").Value = 8.2381033000000006
$ws.Range("J7").Value = 7.8156088999999991
$ws.Range("K7").Value = 669.21583260000011
$ws.Range("L7").Value = 82.082697999999993
$ws.Range("M7").Value = 1318.0059082
$ws.Range("N7").Value = 901.35064690000002
$ws.Range("O7").Value = 6196.2257811999998
$ws.Range("P7").Value = 26143.669140599999
$ws.Range("Q7").Value = -0.45282139999999982
$ws.Range("R7").Value = -0.00020030000000000004

# --- Match the target cell formatting for row 7 ---
# (the newly-inserted row inherited its fill from row 6 above, which is
# already "no fill" for these columns - only the number format needs fixing)
foreach ($col in @("E7", "F7", "G7", "H7", "J7", "Q7")) {
    $ws.Range($col).NumberFormat = "0.00"
}

# Columns with yellow highlight fill, 0.00 number format
foreach ($col in @("D7", "I7", "K7", "L7", "M7", "N7")) {
    $ws.Range($col).NumberFormat = "0.00"
    $ws.Range($col).Interior.Color = 65535
}

# Columns with yellow highlight fill, integer number format
foreach ($col in @("O7", "P7")) {
    $ws.Range($col).NumberFormat = "0"
    $ws.Range($col).Interior.Color = 65535
}

# Mass-balance discrepancy fraction column: no fill, 0.000000 number format
$ws.Range("R7").NumberFormat = "0.000000"

# --- Update selection to match the edited region ---
$ws.Range("K7:P7").Select()
